# Bugfix 21 - accents in report text were not rendered correctly in the
# generated PDF. Update the golden test workbook so that:
#   - the "Data" sheet header row uses accented "cólumn N" labels
#   - the "Summary" sheet contains an accented value ("Número")

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$summarySheet = $wb.Worksheets.Item("Summary")

# 1) Header row A1:L1 on the "Data" sheet: "column N" -> "cólumn N"
for ($col = 1; $col -le 12; $col++) {
    $dataSheet.Cells.Item(1, $col).Value = "cólumn $col"
}

# 2) "Summary" sheet: B4 changes from "Another" to "Número"
$summarySheet.Range("B4").Value = "Número"

$wb.Save()
